# "Small fixes to SLS and ART texts"
#
# The SLS_dict lookup table's KEY_CONTINUE row ("Weiter mit beliebiger
# Taste" / "Press any key to continue") is replaced with text that matches
# the actual "F"/"J" response keys used elsewhere in the questionnaire.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SLS_dict")

# English ("en") value first, German ("de") value second, so the new
# shared-string entries end up appended in that same order.
$ws.Range("C16").Value = "Continue with ""F"" or ""J"""
$ws.Range("B16").Value = "Weiter mit „F“ oder „J“"

# The "de" cell (B16) loses its previous "vertical top" formatting and
# reverts to the default (unstyled) cell style.
$ws.Range("B16").ClearFormats()

# Leave the active selection on the edited cell.
$ws.Range("B16").Select()

# Touch the page setup (paper size / orientation) so it gets persisted.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
